# Adding BMI formula to dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G
    $cell.Formula = "=ROUND((D$row/1.88)/1.88,2)"
    $cell.Style = "Normal"
}
